# Generate Report for Handback
# Updates the handback-status workbook: the two source files that were
# processed in this run were renamed/regenerated (new GUID-based names and
# new timestamps), so every cell/hyperlink that referenced the old names
# needs to move to the new ones.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Mapping of OLD display text -> NEW display text for every hyperlink in
# the workbook. The hyperlink's target Address (rId) is left untouched -
# only the visible text changes, same as the canonical edit.
# ---------------------------------------------------------------------
function Update-SheetHyperlinks {
    param($ws)

    foreach ($h in $ws.Hyperlinks) {
        $old = $h.TextToDisplay()

        if ($old -eq "e2e\0df0035a-b1b3-4f8c-86a4-cbe0abc559e8.md") {
            $h.TextToDisplay = "e2e\57340785-7b4f-457d-9508-536245ce2b9a.md"
        }
        elseif ($old -eq "e2e\574d7319-3b13-4ec4-8bc9-99d88c62241f.md") {
            $h.TextToDisplay = "e2e\ffff9f9a5b3d-23e1-4624-8221-cc52b40e66b9.md"
        }
        elseif ($old -eq "0df0035a-b1b3-4f8c-86a4-cbe0abc559e8.md") {
            $h.TextToDisplay = "57340785-7b4f-457d-9508-536245ce2b9a.md"
        }
        elseif ($old -eq "574d7319-3b13-4ec4-8bc9-99d88c62241f.md") {
            $h.TextToDisplay = "ffff9f9a5b3d-23e1-4624-8221-cc52b40e66b9.md"
        }
    }
}

# ----------------------------- Overview -------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
Update-SheetHyperlinks $wsOverview

$wsOverview.Range("A2").Value = "57340785-7b4f-457d-9508-536245ce2b9a.md"
$wsOverview.Range("B2").Value = "e2e\57340785-7b4f-457d-9508-536245ce2b9a.md"
$wsOverview.Range("G2").Value = "2016-08-23 03:01:09"

$wsOverview.Range("A3").Value = "ffff9f9a5b3d-23e1-4624-8221-cc52b40e66b9.md"
$wsOverview.Range("B3").Value = "e2e\ffff9f9a5b3d-23e1-4624-8221-cc52b40e66b9.md"
$wsOverview.Range("G3").Value = "2016-08-23 03:01:09"

# ------------------------------- zh-cn ---------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
Update-SheetHyperlinks $wsZhCn

$wsZhCn.Range("A2").Value = "57340785-7b4f-457d-9508-536245ce2b9a.md"
$wsZhCn.Range("G2").Value = "57340785-7b4f-457d-9508-536245ce2b9a.c3400fcc7e934cd93d4997cf0adb20e65eadf6a5.zh-cn.xlf"
$wsZhCn.Range("H2").Value = "2016-08-23 03:00:59"
$wsZhCn.Range("I2").Value = "57340785-7b4f-457d-9508-536245ce2b9a.md"
$wsZhCn.Range("J2").Value = "57340785-7b4f-457d-9508-536245ce2b9a.c3400fcc7e934cd93d4997cf0adb20e65eadf6a5.zh-cn.xlf"
$wsZhCn.Range("K2").Value = "2016-08-23 03:01:28"

$wsZhCn.Range("A3").Value = "ffff9f9a5b3d-23e1-4624-8221-cc52b40e66b9.md"
$wsZhCn.Range("G3").Value = "57340785-7b4f-457d-9508-536245ce2b9a.c3400fcc7e934cd93d4997cf0adb20e65eadf6a5.zh-cn.xlf"
$wsZhCn.Range("H3").Value = "2016-08-23 03:00:59"
$wsZhCn.Range("I3").Value = "ffff9f9a5b3d-23e1-4624-8221-cc52b40e66b9.md"
$wsZhCn.Range("J3").Value = "57340785-7b4f-457d-9508-536245ce2b9a.c3400fcc7e934cd93d4997cf0adb20e65eadf6a5.zh-cn.xlf"
$wsZhCn.Range("K3").Value = "2016-08-23 03:01:28"

# ------------------------------- de-de ---------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
Update-SheetHyperlinks $wsDeDe

$wsDeDe.Range("A2").Value = "57340785-7b4f-457d-9508-536245ce2b9a.md"
$wsDeDe.Range("G2").Value = "57340785-7b4f-457d-9508-536245ce2b9a.c3400fcc7e934cd93d4997cf0adb20e65eadf6a5.de-de.xlf"
$wsDeDe.Range("H2").Value = "2016-08-23 03:01:09"
$wsDeDe.Range("I2").Value = "57340785-7b4f-457d-9508-536245ce2b9a.md"
$wsDeDe.Range("J2").Value = "57340785-7b4f-457d-9508-536245ce2b9a.c3400fcc7e934cd93d4997cf0adb20e65eadf6a5.de-de.xlf"
$wsDeDe.Range("K2").Value = "2016-08-23 03:01:35"

$wsDeDe.Range("A3").Value = "ffff9f9a5b3d-23e1-4624-8221-cc52b40e66b9.md"
$wsDeDe.Range("G3").Value = "57340785-7b4f-457d-9508-536245ce2b9a.c3400fcc7e934cd93d4997cf0adb20e65eadf6a5.de-de.xlf"
$wsDeDe.Range("H3").Value = "2016-08-23 03:01:09"
$wsDeDe.Range("I3").Value = "ffff9f9a5b3d-23e1-4624-8221-cc52b40e66b9.md"
$wsDeDe.Range("J3").Value = "57340785-7b4f-457d-9508-536245ce2b9a.c3400fcc7e934cd93d4997cf0adb20e65eadf6a5.de-de.xlf"
$wsDeDe.Range("K3").Value = "2016-08-23 03:01:35"
